$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 113. This shifts the existing rows
# 113-145 down to 114-146 (their content/formatting is preserved by Excel).
$ws.Rows.Item(113).Insert()

# Populate the newly inserted row 113 with the new price observation.
$ws.Cells.Item(113, 1).Value = 4
$ws.Cells.Item(113, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(113, 3).Value = "Los Lagos"
$ws.Cells.Item(113, 4).Value = 44876
$ws.Cells.Item(113, 5).Value = 10
$ws.Cells.Item(113, 6).Value = 100112022
$ws.Cells.Item(113, 7).Value = "Arveja Verde"
$ws.Cells.Item(113, 8).Value = "Perfection"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 80
$ws.Cells.Item(113, 11).Value = 25000
$ws.Cells.Item(113, 12).Value = 25000
$ws.Cells.Item(113, 13).Value = 25000
$ws.Cells.Item(113, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(113, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(113, 16).Value = 1000
$ws.Cells.Item(113, 17).Value = 25
$ws.Cells.Item(113, 18).Value = "Hortaliza"
